$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = '[Object Detection] YOLO Define Optimal Anchor Box :: YOLO v5, YOLO v6 autoanchor'
$ws.Range("D16").Value = '[GAN] Stitch it in Time: GAN-Based Facial Editing of Real Videos'
$ws.Range("D19").Value = '아기 있는 집에는 에몬스홈 그란데 가죽소파 4인'
$ws.Range("D20").Value = '[책][리뷰] Developer Relations'
$ws.Range("D23").Value = '[arxiv.org에 올라온 논문을 pdf가 아닌 잘 정돈된 1 column의 텍스트 페이지로 볼 수 있는 방법] 주소 창에서 arxiv의 x를 숫자 5로 바꾸면 pdf가 아닌 텍스트 페이지로 나옴'
$ws.Range("E23").Value = 'https://theonly1.tistory.com/3064'
$ws.Range("D28").Value = 'Mobile Manipulator 101 :: Dual Trajectory & task priority(2)'
$ws.Range("E28").Value = 'https://ropiens.tistory.com/201'
$ws.Range("D32").Value = 'Hadoop Ecosystem 하둡 에코시스템 간단정리'
$ws.Range("D39").Value = 'Visualize your data with Facets'
$ws.Range("D42").Value = '위경도 도분초, 십진법 값별 거리 차이'
$ws.Range("D43").Value = 'np.random.shuffle 과 np.random.permutation 정리'
$ws.Range("D44").Value = 'Object Detection Algorithm (Efficientdet)'
$ws.Range("D45").Value = 'RNN Auto-Encoder (RAE)'
$ws.Range("D46").Value = '[국립암센터] 2022년 10월, 생물정보학(Bioinformatics 채용), 연구소 정규직 연구직'
$ws.Range("E46").Value = 'https://bioinformaticsandme.tistory.com/492'
$ws.Range("D47").Value = 'Pseudo Labeling, TTA(Test Time Augmentation) 기법'
$ws.Range("D51").Value = '마크다운으로 깔끔하게 글을 쓸 수 있는 블로그 플랫폼, velog'
$ws.Range("E51").Value = 'https://bskyvision.com/entry/%EB%A7%88%ED%81%AC%EB%8B%A4%EC%9A%B4%EC%9C%BC%EB%A1%9C-%EA%B9%94%EB%81%94%ED%95%98%EA%B2%8C-%EA%B8%80%EC%9D%84-%EC%93%B8-%EC%88%98-%EC%9E%88%EB%8A%94-%EB%B8%94%EB%A1%9C%EA%B7%B8-%ED%94%8C%EB%9E%AB%ED%8F%BC-velog'
